# gradient descent with momentum (incl. bias issue resolved)
#
# Adds 4 new benchmark rows (19-22, dated 2018.04.30) to Sheet1, plus the
# accompanying comments, and a couple of small view/format tweaks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Helper: write a date-look-alike string ("2018.04.30") as plain TEXT.
# A bare .Value assignment of that literal gets auto-parsed into a date
# serial by the recalc engine, so we park it (apostrophe-prefixed, which
# forces text) in a scratch cell in the same column, copy it, and paste
# *values only* into the destination -- PasteSpecial(xlPasteValues=-4163)
# carries the resolved text over without re-triggering date detection
# and without dragging along the scratch cell's quote-prefix formatting.
# ---------------------------------------------------------------------
function Set-DateLabel($target) {
    $ws.Range("B30").Value = "'2018.04.30"
    $ws.Range("B30").Copy()
    $ws.Range($target).PasteSpecial(-4163)
}

# Row 19 --------------------------------------------------------------
Set-DateLabel "B19"
$ws.Range("C19").Value = 2.6
$ws.Range("D19").Value = 0.75608305224674899
$ws.Range("E19").Value = 0.80400000000000005
$ws.Range("F19").Value = 1000
$ws.Range("G19").Value = 0.001
$ws.Range("H19").Value = "tanh*,softmax"
$ws.Range("I19").Value = "forRelu"
$ws.Range("J19").Value = 2
$ws.Range("K19").Value = "'100, 10"

# Row 20 --------------------------------------------------------------
Set-DateLabel "B20"
$ws.Range("C20").Value = 2.9
$ws.Range("D20").Value = 0.69801583265231404
$ws.Range("E20").Value = 0.83399999999999996
$ws.Range("F20").Value = 1000
$ws.Range("G20").Value = 0.001
$ws.Range("H20").Value = "relu*,softmax"
$ws.Range("I20").Value = "forRelu"
$ws.Range("J20").Value = 2
$ws.Range("K20").Value = "'100, 10"
$ws.Range("L20").Value = "relu slightly better than tanh"

# Row 21 --------------------------------------------------------------
Set-DateLabel "B21"
$ws.Range("C21").Value = 2.5
$ws.Range("D21").Value = 0.85885793959385504
$ws.Range("E21").Value = 0.82899999999999996
$ws.Range("F21").Value = 1000
$ws.Range("G21").Value = 0.001
# H21 is a two-run rich string: bold "tanh*" + regular ",softmax"
$ws.Range("H21").Value = "tanh*,softmax"
$ws.Range("H21").Characters(1, 5).Font.Bold = $true
$ws.Range("I21").Value = "xavier"
# Match the bold "xavier" styling used for the same label in I10.
$ws.Range("I10").Copy()
$ws.Range("I21").PasteSpecial(-4122)
$ws.Range("J21").Value = 2
$ws.Range("K21").Value = "'100, 10"
$ws.Range("L21").Value = "xavier did not help much, but maybe because alpha too high?"

# Row 22 --------------------------------------------------------------
Set-DateLabel "B22"
$ws.Range("C22").Value = 3.1
$ws.Range("D22").Value = 0.83597837807029796
$ws.Range("E22").Value = "?"
$ws.Range("F22").Value = 1000
$ws.Range("G22").Value = 0.001
$ws.Range("H22").Value = "relu*,softmax"
$ws.Range("I22").Value = "forRelu"
$ws.Range("J22").Value = 2
$ws.Range("K22").Value = "'100, 10"
$ws.Range("L22").Value = "MINIBATCHES implemented (up to here: trained on 1K data set, costs jumped between ,75 and 1,00 in last iterations"

# Clean up the scratch cell used by Set-DateLabel.
$ws.Range("B30").Clear()

# Column K is a touch narrower now (best-fit to the "100, 10" labels).
$ws.Columns("K").ColumnWidth = 10.5

# Selection follows the last-used row, same as the live editing session.
$ws.Range("E23").Select() | Out-Null
